$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 24394010
$ws.Range("I98").Value = 3372.9167
$ws.Range("J98").Value = 200006600
$ws.Range("K98").Value = 3372.9167
$ws.Range("L98").Value = 200006600
$ws.Range("M98").Value = -1874.9167
$ws.Range("N98").Value = -200009596

$ws.Range("H122").Value = 24394010
$ws.Range("I122").Value = 3372.9167
$ws.Range("J122").Value = 200006600
$ws.Range("K122").Value = 10118.7501
$ws.Range("L122").Value = 600019800
$ws.Range("M122").Value = -7668.750100000001
$ws.Range("N122").Value = -600024700

$ws.Range("H132").Value = 3083.6223
$ws.Range("I132").Value = 3223.8157
$ws.Range("J132").Value = 2322.5715
$ws.Range("K132").Value = 9671.447100000001
$ws.Range("L132").Value = 6967.7145
$ws.Range("M132").Value = -7141.447100000001
$ws.Range("N132").Value = -12027.7145

$ws.Range("H137").Value = 1803.9584
$ws.Range("I137").Value = 1849.8
$ws.Range("J137").Value = 1574.75
$ws.Range("K137").Value = 5549.4
$ws.Range("L137").Value = 4724.25
$ws.Range("M137").Value = -2999.4
$ws.Range("N137").Value = -9824.25

$ws.Range("H138").Value = 4227.6445
$ws.Range("I138").Value = 1485.4865
$ws.Range("J138").Value = 6829.1797
$ws.Range("K138").Value = 4456.4595
$ws.Range("L138").Value = 20487.5391
$ws.Range("M138").Value = 683.5405000000001
$ws.Range("N138").Value = -30767.5391

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 351272.72
$ws.Range("I32").Value = 2406.0532
$ws.Range("K32").Value = 2406.0532
$ws.Range("M32").Value = -2119.0532

$ws.Range("H61").Value = 1220541.2
$ws.Range("I61").Value = 1725000
$ws.Range("J61").Value = 1432.5
$ws.Range("K61").Value = 1725000
$ws.Range("L61").Value = 1432.5
$ws.Range("M61").Value = -1724788
$ws.Range("N61").Value = -1856.5

$ws.Range("H136").Value = 1220541.2
$ws.Range("I136").Value = 1725000
$ws.Range("J136").Value = 1432.5
$ws.Range("K136").Value = 5175000
$ws.Range("L136").Value = 4297.5
$ws.Range("M136").Value = -5172450
$ws.Range("N136").Value = -9397.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2293.5134
$ws.Range("I31").Value = 1745.6296
$ws.Range("J31").Value = 3772.8
$ws.Range("K31").Value = 1745.6296
$ws.Range("L31").Value = 3772.8
$ws.Range("M31").Value = -1450.6296
$ws.Range("N31").Value = -4362.8

$ws.Range("H34").Value = 2293.5134
$ws.Range("I34").Value = 1745.6296
$ws.Range("J34").Value = 3772.8
$ws.Range("K34").Value = 1745.6296
$ws.Range("L34").Value = 3772.8
$ws.Range("M34").Value = -1543.6296
$ws.Range("N34").Value = -4176.8

$ws.Range("H107").Value = 572.92
$ws.Range("I107").Value = 457.92856
$ws.Range("J107").Value = 719.2727
$ws.Range("K107").Value = 457.92856
$ws.Range("L107").Value = 719.2727
$ws.Range("M107").Value = 1462.07144
$ws.Range("N107").Value = -4559.2727

$ws.Range("H132").Value = 17056.111
$ws.Range("I132").Value = 861.70215
$ws.Range("J132").Value = 64627.188
$ws.Range("K132").Value = 2585.10645
$ws.Range("L132").Value = 193881.564
$ws.Range("M132").Value = -55.10644999999977
$ws.Range("N132").Value = -198941.564

$ws.Range("H134").Value = 296025.5
$ws.Range("I134").Value = 323749.6
$ws.Range("J134").Value = 9543.333000000001
$ws.Range("K134").Value = 971248.7999999999
$ws.Range("L134").Value = 28629.999
$ws.Range("M134").Value = -968713.7999999999
$ws.Range("N134").Value = -33699.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 446.85715
$ws.Range("I23").Value = 156.33333
$ws.Range("J23").Value = 664.75
$ws.Range("K23").Value = 468.99999
$ws.Range("L23").Value = 1994.25
$ws.Range("M23").Value = -233.99999
$ws.Range("N23").Value = -2464.25

$ws.Range("H70").Value = 1343.6666
$ws.Range("I70").Value = 687.3333
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 2061.9999
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -1746.9999
$ws.Range("N70").Value = -6630

$ws.Range("H73").Value = 1343.6666
$ws.Range("I73").Value = 687.3333
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 2061.9999
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -969.9998999999998
$ws.Range("N73").Value = -8184

$ws.Range("H75").Value = 1776.25
$ws.Range("J75").Value = 1776.25
$ws.Range("L75").Value = 5328.75
$ws.Range("N75").Value = -7324.75

$ws.Range("H78").Value = 1776.25
$ws.Range("J78").Value = 1776.25
$ws.Range("L78").Value = 15986.25
$ws.Range("N78").Value = -25970.25

$ws.Range("H113").Value = 838.1486
$ws.Range("I113").Value = 616.5
$ws.Range("J113").Value = 973.06525
$ws.Range("K113").Value = 1849.5
$ws.Range("L113").Value = 2919.19575
$ws.Range("M113").Value = 320.5
$ws.Range("N113").Value = -7259.19575

$ws.Range("H114").Value = 561680.6
$ws.Range("I114").Value = 10409.8
$ws.Range("J114").Value = 1250769.1
$ws.Range("K114").Value = 31229.4
$ws.Range("L114").Value = 3752307.3
$ws.Range("M114").Value = -27975.4
$ws.Range("N114").Value = -3758815.3

$ws.Range("H131").Value = 62502350
$ws.Range("I131").Value = 1725
$ws.Range("J131").Value = 125002980
$ws.Range("K131").Value = 5175
$ws.Range("L131").Value = 375008940
$ws.Range("M131").Value = -135
$ws.Range("N131").Value = -375019020

$ws.Range("H132").Value = 2392.9
$ws.Range("I132").Value = 1650
$ws.Range("J132").Value = 2711.2856
$ws.Range("K132").Value = 14850
$ws.Range("L132").Value = 24401.5704
$ws.Range("M132").Value = -12320
$ws.Range("N132").Value = -29461.5704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4356.1514
$ws.Range("I132").Value = 3168.4375
$ws.Range("J132").Value = 7523.3887
$ws.Range("K132").Value = 9505.3125
$ws.Range("L132").Value = 22570.1661
$ws.Range("M132").Value = -6975.3125
$ws.Range("N132").Value = -27630.1661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H132").Value = 25658.818
$ws.Range("I132").Value = 40734.035
$ws.Range("J132").Value = 1715.8235
$ws.Range("K132").Value = 122202.105
$ws.Range("L132").Value = 5147.470499999999
$ws.Range("M132").Value = -119672.105
$ws.Range("N132").Value = -10207.4705

$ws.Range("H136").Value = 5353.362
$ws.Range("I136").Value = 5468.909
$ws.Range("J136").Value = 5081
$ws.Range("K136").Value = 16406.727
$ws.Range("L136").Value = 15243
$ws.Range("M136").Value = -13856.727
$ws.Range("N136").Value = -20343

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 30056352
$ws.Range("I132").Value = 40646320
$ws.Range("J132").Value = 2698939
$ws.Range("K132").Value = 121938960
$ws.Range("L132").Value = 8096817
$ws.Range("M132").Value = -121936430
$ws.Range("N132").Value = -8101877
